$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.348259806632996
$ws.Range("B1").Value = 2.86980938911438
$ws.Range("C1").Value = 3.989838123321533
$ws.Range("D1").Value = 0.9707400798797607
$ws.Range("E1").Value = 0.6288301348686218
